$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells I1 and J1, copying the format of the existing H1 header cell
# (bold font, thin border, centered alignment) then set their text.
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("I1").Value = "I0"

$ws.Range("H1").Copy($ws.Range("J1"))
$ws.Range("J1").Value = "IF"

# Data for columns I (I0) and J (IF), for rows 2-63
$data = @(
    @(4, 5),
    @(7, 8),
    @(6, 6),
    @(8, 9),
    @(9, 9),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(7, 8),
    @(8, 8),
    @(8, 8),
    @(7, 7),
    @(9, 9),
    @(8, 9),
    @(7, 8),
    @(9, 9),
    @(6, 7),
    @(9, 9),
    @(8, 8),
    @(9, 9),
    @(5, 5),
    @(9, 9),
    @(8, 8),
    @(8, 8),
    @(5, 5),
    @(9, 9),
    @(6, 6),
    @(7, 7),
    @(9, 9),
    @(6, 6),
    @(3, 3),
    @(8, 8),
    @(6, 7),
    @(7, 7),
    @(5, 5),
    @(7, 7),
    @(7, 7),
    @(7, 7),
    @(7, 7),
    @(8, 8),
    @(6, 7),
    @(5, 5),
    @(6, 6),
    @(8, 8),
    @(9, 9),
    @(4, 4),
    @(9, 9),
    @(8, 8),
    @(4, 5),
    @(8, 8),
    @(9, 9),
    @(7, 7),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(7, 7),
    @(7, 7),
    @(5, 5),
    @(8, 8),
    @(8, 8),
    @(7, 7),
    @(8, 8)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 9).Value = $row[0]
    $ws.Cells.Item($r, 10).Value = $row[1]
    $r = $r + 1
}
